# Stock-count correction pass.
#
# Per item row: F = on-hand quantity, G = D * F (stock value at cost).
# Each company block ends in a "Sub Total:" row whose B cell is the sum
# of G across that block; B619/B620 are the grand total across every
# company subtotal. A handful of same-product row pairs (e.g. rows
# 127/128) had their batch code / rate / qty / value (B,E,F,G) swapped
# between the two rows - serial number (A) and product name (C) stay put.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F9").Value = 19
$ws.Range("G9").Value = 561.83
$ws.Range("B10").Value = 33787.18
$ws.Range("F55").Value = 2
$ws.Range("G55").Value = 501.58
$ws.Range("F56").Value = 15
$ws.Range("G56").Value = 3137.7
$ws.Range("F68").Value = 68
$ws.Range("G68").Value = 7828.16
$ws.Range("F80").Value = 16
$ws.Range("G80").Value = 3937.12
$ws.Range("F81").Value = 29
$ws.Range("G81").Value = 886.8200000000001
$ws.Range("F83").Value = 123
$ws.Range("G83").Value = 18532.41
$ws.Range("F87").Value = 17
$ws.Range("G87").Value = 5232.09
$ws.Range("B90").Value = 212563.27
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 0
$ws.Range("B100").Value = 2115.59
$ws.Range("F115").Value = 241
$ws.Range("G115").Value = 23331.21
$ws.Range("B117").Value = 17498.67
$ws.Range("B127").Value = 57552
$ws.Range("E127").Value = 136.86
$ws.Range("F127").Value = -5
$ws.Range("G127").Value = -603.45
$ws.Range("B128").Value = 64329
$ws.Range("E128").Value = 128.32
$ws.Range("F128").Value = 2
$ws.Range("G128").Value = 241.38
$ws.Range("F146").Value = 44
$ws.Range("G146").Value = 3704.36
$ws.Range("B147").Value = 25167.41
$ws.Range("F149").Value = 263
$ws.Range("G149").Value = 17042.4
$ws.Range("F151").Value = 108
$ws.Range("G151").Value = 9383.040000000001
$ws.Range("F152").Value = 78
$ws.Range("G152").Value = 6886.62
$ws.Range("B156").Value = 37752.72
$ws.Range("B192").Value = 64973
$ws.Range("E192").Value = 35.4
$ws.Range("F192").Value = 2
$ws.Range("G192").Value = 66.59999999999999
$ws.Range("B193").Value = 48706
$ws.Range("E193").Value = 39.8
$ws.Range("F193").Value = -144
$ws.Range("G193").Value = -4795.2
$ws.Range("F196").Value = 3
$ws.Range("G196").Value = 343.68
$ws.Range("F201").Value = 0
$ws.Range("G201").Value = 0
$ws.Range("F203").Value = 81
$ws.Range("G203").Value = 1632.96
$ws.Range("F206").Value = 0
$ws.Range("G206").Value = 0
$ws.Range("F213").Value = 10
$ws.Range("G213").Value = 856.8
$ws.Range("F214").Value = 56
$ws.Range("G214").Value = 4911.2
$ws.Range("B216").Value = 57296.19
$ws.Range("F233").Value = 128
$ws.Range("G233").Value = 6097.92
$ws.Range("F248").Value = 5
$ws.Range("G248").Value = 295.65
$ws.Range("F255").Value = 631
$ws.Range("G255").Value = 108109.23
$ws.Range("F256").Value = 307
$ws.Range("G256").Value = 46409.19
$ws.Range("B260").Value = 223498.26
$ws.Range("F270").Value = 57
$ws.Range("G270").Value = 1837.68
$ws.Range("B275").Value = 10077.94
$ws.Range("F277").Value = 7
$ws.Range("G277").Value = 148.75
$ws.Range("F282").Value = 20
$ws.Range("G282").Value = 1074
$ws.Range("F285").Value = 38
$ws.Range("G285").Value = 1061.34
$ws.Range("F292").Value = 52
$ws.Range("G292").Value = 4330.04
$ws.Range("F294").Value = 56
$ws.Range("G294").Value = 3996.16
$ws.Range("F295").Value = 9
$ws.Range("G295").Value = 933.21
$ws.Range("F296").Value = 108
$ws.Range("G296").Value = 2289.6
$ws.Range("F302").Value = 84
$ws.Range("G302").Value = 17714.76
$ws.Range("B304").Value = 207412.34
$ws.Range("F307").Value = 3
$ws.Range("G307").Value = 462.45
$ws.Range("B309").Value = 2551.91
$ws.Range("F328").Value = 74
$ws.Range("G328").Value = 2753.54
$ws.Range("B330").Value = 33494.02
$ws.Range("F341").Value = 9
$ws.Range("G341").Value = 458.55
$ws.Range("F345").Value = 96
$ws.Range("G345").Value = 5895.36
$ws.Range("B346").Value = 31777.05
$ws.Range("F350").Value = 65
$ws.Range("G350").Value = 4987.45
$ws.Range("B358").Value = 39492.45
$ws.Range("B366").Value = 65066
$ws.Range("E366").Value = 13.61
$ws.Range("F366").Value = 90
$ws.Range("G366").Value = 1152.9
$ws.Range("B367").Value = 53263
$ws.Range("E367").Value = 15.29
$ws.Range("F367").Value = -309
$ws.Range("G367").Value = -3958.29
$ws.Range("B375").Value = 45718
$ws.Range("E375").Value = 19.38
$ws.Range("F375").Value = -294
$ws.Range("G375").Value = -4768.68
$ws.Range("B376").Value = 64927
$ws.Range("E376").Value = 17.26
$ws.Range("F376").Value = 106
$ws.Range("G376").Value = 1719.32
$ws.Range("B380").Value = 64925
$ws.Range("E380").Value = 13.97
$ws.Range("F380").Value = 111
$ws.Range("G380").Value = 1459.65
$ws.Range("B381").Value = 45709
$ws.Range("E381").Value = 15.69
$ws.Range("F381").Value = -300
$ws.Range("G381").Value = -3945
$ws.Range("F390").Value = 18
$ws.Range("G390").Value = 1109.7
$ws.Range("F394").Value = 0
$ws.Range("G394").Value = 0
$ws.Range("B395").Value = 2126.81
$ws.Range("F429").Value = 26
$ws.Range("G429").Value = 176.8
$ws.Range("F430").Value = 24
$ws.Range("G430").Value = 309.36
$ws.Range("B435").Value = 2336.06
$ws.Range("B442").Value = 64810
$ws.Range("E442").Value = 291.22
$ws.Range("F442").Value = 5
$ws.Range("G442").Value = 1369.6
$ws.Range("B443").Value = 53319
$ws.Range("E443").Value = 310.64
$ws.Range("F443").Value = -6
$ws.Range("G443").Value = -1643.52
$ws.Range("F477").Value = 20
$ws.Range("G477").Value = 906.8
$ws.Range("B478").Value = 906.8
$ws.Range("F482").Value = 47
$ws.Range("G482").Value = 2785.69
$ws.Range("F485").Value = 32
$ws.Range("G485").Value = 5615.04
$ws.Range("B488").Value = 34055.14
$ws.Range("F490").Value = 12
$ws.Range("G490").Value = 1566.6
$ws.Range("B493").Value = 15950.76
$ws.Range("F542").Value = 56
$ws.Range("G542").Value = 7253.68
$ws.Range("F544").Value = 1
$ws.Range("G544").Value = 794.48
$ws.Range("B547").Value = 25910.12
$ws.Range("F551").Value = 34
$ws.Range("G551").Value = 4866.42
$ws.Range("F553").Value = 1
$ws.Range("G553").Value = 75.68000000000001
$ws.Range("F558").Value = 83
$ws.Range("G558").Value = 11205.83
$ws.Range("B560").Value = 26368.36
$ws.Range("F575").Value = 17
$ws.Range("G575").Value = 562.02
$ws.Range("F578").Value = 107
$ws.Range("G578").Value = 5338.23
$ws.Range("F581").Value = 37
$ws.Range("G581").Value = 8946.6
$ws.Range("F582").Value = 67
$ws.Range("G582").Value = 3818.33
$ws.Range("B583").Value = 35271.41
$ws.Range("F599").Value = 2375
$ws.Range("G599").Value = 387386.25
$ws.Range("F601").Value = 501
$ws.Range("G601").Value = 141717.87
$ws.Range("F602").Value = 382
$ws.Range("G602").Value = 55256.3
$ws.Range("B606").Value = 585208.47
$ws.Range("F613").Value = 164
$ws.Range("G613").Value = 26102.24
$ws.Range("F617").Value = 20
$ws.Range("G617").Value = 790.6
$ws.Range("B618").Value = 51709.16
$ws.Range("B619").Value = 2169384.7
$ws.Range("B620").Value = 2169384.7
